$d = $word.ActiveDocument

$d.Content.Find.Execute("22-15=", $true, $false, $false, $false, $false, $true, 1, $false, "58-24=", 2)
$d.Content.Find.Execute("73-62=", $true, $false, $false, $false, $false, $true, 1, $false, "57-11=", 2)
$d.Content.Find.Execute("40+26=", $true, $false, $false, $false, $false, $true, 1, $false, "52-3=", 2)
$d.Content.Find.Execute("34+45=", $true, $false, $false, $false, $false, $true, 1, $false, "8+7=", 2)
$d.Content.Find.Execute("15+10=", $true, $false, $false, $false, $false, $true, 1, $false, "27+47=", 2)
$d.Content.Find.Execute("2+1=", $true, $false, $false, $false, $false, $true, 1, $false, "65+23=", 2)
$d.Content.Find.Execute("76-6=", $true, $false, $false, $false, $false, $true, 1, $false, "80+6=", 2)
$d.Content.Find.Execute("41-39=", $true, $false, $false, $false, $false, $true, 1, $false, "80-47=", 2)
$d.Content.Find.Execute("43-15=", $true, $false, $false, $false, $false, $true, 1, $false, "81+9=", 2)
$d.Content.Find.Execute("30+42=", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=", 2)
$d.Content.Find.Execute("51-9=", $true, $false, $false, $false, $false, $true, 1, $false, "33+51=", 2)
$d.Content.Find.Execute("18+42=", $true, $false, $false, $false, $false, $true, 1, $false, "11+23=", 2)
$d.Content.Find.Execute("93-36=", $true, $false, $false, $false, $false, $true, 1, $false, "55+12=", 2)
$d.Content.Find.Execute("58-49=", $true, $false, $false, $false, $false, $true, 1, $false, "6+13=", 2)
$d.Content.Find.Execute("95-61=", $true, $false, $false, $false, $false, $true, 1, $false, "93-23=", 2)
$d.Content.Find.Execute("85+13=", $true, $false, $false, $false, $false, $true, 1, $false, "49+44=", 2)
$d.Content.Find.Execute("23+6=", $true, $false, $false, $false, $false, $true, 1, $false, "22+14=", 2)
$d.Content.Find.Execute("94-72=", $true, $false, $false, $false, $false, $true, 1, $false, "49-23=", 2)
$d.Content.Find.Execute("59+14=", $true, $false, $false, $false, $false, $true, 1, $false, "38-22=", 2)
$d.Content.Find.Execute("78-2=", $true, $false, $false, $false, $false, $true, 1, $false, "47+52=", 2)
$d.Content.Find.Execute("3+29=", $true, $false, $false, $false, $false, $true, 1, $false, "39+50=", 2)
$d.Content.Find.Execute("53-32=", $true, $false, $false, $false, $false, $true, 1, $false, "41-37=", 2)
$d.Content.Find.Execute("52+40=", $true, $false, $false, $false, $false, $true, 1, $false, "3+82=", 2)
$d.Content.Find.Execute("96-3=", $true, $false, $false, $false, $false, $true, 1, $false, "12+7=", 2)
$d.Content.Find.Execute("82+14=", $true, $false, $false, $false, $false, $true, 1, $false, "21+50=", 2)
$d.Content.Find.Execute("24+37=", $true, $false, $false, $false, $false, $true, 1, $false, "69-43=", 2)
$d.Content.Find.Execute("25+73=", $true, $false, $false, $false, $false, $true, 1, $false, "75+19=", 2)
$d.Content.Find.Execute("3+52=", $true, $false, $false, $false, $false, $true, 1, $false, "46+14=", 2)
$d.Content.Find.Execute("70+7=", $true, $false, $false, $false, $false, $true, 1, $false, "35+31=", 2)
$d.Content.Find.Execute("67+31=", $true, $false, $false, $false, $false, $true, 1, $false, "73-57=", 2)
$d.Content.Find.Execute("51-33=", $true, $false, $false, $false, $false, $true, 1, $false, "10+62=", 2)
$d.Content.Find.Execute("30-13=", $true, $false, $false, $false, $false, $true, 1, $false, "13-4=", 2)
$d.Content.Find.Execute("36-26=", $true, $false, $false, $false, $false, $true, 1, $false, "57-10=", 2)
$d.Content.Find.Execute("43-20=", $true, $false, $false, $false, $false, $true, 1, $false, "56-26=", 2)
$d.Content.Find.Execute("15+24=", $true, $false, $false, $false, $false, $true, 1, $false, "72+1=", 2)
$d.Content.Find.Execute("49+19=", $true, $false, $false, $false, $false, $true, 1, $false, "33+28=", 2)
$d.Content.Find.Execute("65-51=", $true, $false, $false, $false, $false, $true, 1, $false, "72+5=", 2)
$d.Content.Find.Execute("49+15=", $true, $false, $false, $false, $false, $true, 1, $false, "89-43=", 2)
$d.Content.Find.Execute("0+93=", $true, $false, $false, $false, $false, $true, 1, $false, "6+83=", 2)
$d.Content.Find.Execute("7-6=", $true, $false, $false, $false, $false, $true, 1, $false, "48+33=", 2)
$d.Content.Find.Execute("60-23=", $true, $false, $false, $false, $false, $true, 1, $false, "55-41=", 2)
$d.Content.Find.Execute("96-89=", $true, $false, $false, $false, $false, $true, 1, $false, "69-12=", 2)
$d.Content.Find.Execute("7+67=", $true, $false, $false, $false, $false, $true, 1, $false, "2+70=", 2)
$d.Content.Find.Execute("63-0=", $true, $false, $false, $false, $false, $true, 1, $false, "44+12=", 2)
$d.Content.Find.Execute("35-2=", $true, $false, $false, $false, $false, $true, 1, $false, "87-65=", 2)
$d.Content.Find.Execute("58-19=", $true, $false, $false, $false, $false, $true, 1, $false, "13+50=", 2)
$d.Content.Find.Execute("55-7=", $true, $false, $false, $false, $false, $true, 1, $false, "47+5=", 2)
$d.Content.Find.Execute("80+11=", $true, $false, $false, $false, $false, $true, 1, $false, "32+38=", 2)
$d.Content.Find.Execute("59-0=", $true, $false, $false, $false, $false, $true, 1, $false, "36+55=", 2)
$d.Content.Find.Execute("6+43=", $true, $false, $false, $false, $false, $true, 1, $false, "95-43=", 2)
$d.Content.Find.Execute("1+11=", $true, $false, $false, $false, $false, $true, 1, $false, "40+54=", 2)
$d.Content.Find.Execute("16-12=", $true, $false, $false, $false, $false, $true, 1, $false, "48+8=", 2)
$d.Content.Find.Execute("31+22=", $true, $false, $false, $false, $false, $true, 1, $false, "38-28=", 2)
$d.Content.Find.Execute("72-4=", $true, $false, $false, $false, $false, $true, 1, $false, "20+66=", 2)
$d.Content.Find.Execute("36+39=", $true, $false, $false, $false, $false, $true, 1, $false, "94-37=", 2)
$d.Content.Find.Execute("14-10=", $true, $false, $false, $false, $false, $true, 1, $false, "81-69=", 2)
$d.Content.Find.Execute("90-46=", $true, $false, $false, $false, $false, $true, 1, $false, "79+14=", 2)
$d.Content.Find.Execute("81-11=", $true, $false, $false, $false, $false, $true, 1, $false, "55+19=", 2)
$d.Content.Find.Execute("38+21=", $true, $false, $false, $false, $false, $true, 1, $false, "34-0=", 2)
$d.Content.Find.Execute("26-24=", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=", 2)
$d.Content.Find.Execute("33-26=", $true, $false, $false, $false, $false, $true, 1, $false, "43+8=", 2)
$d.Content.Find.Execute("11-6=", $true, $false, $false, $false, $false, $true, 1, $false, "46-29=", 2)
$d.Content.Find.Execute("84-65=", $true, $false, $false, $false, $false, $true, 1, $false, "87-25=", 2)
$d.Content.Find.Execute("48-30=", $true, $false, $false, $false, $false, $true, 1, $false, "61-18=", 2)
$d.Content.Find.Execute("7-4=", $true, $false, $false, $false, $false, $true, 1, $false, "63+19=", 2)
$d.Content.Find.Execute("42+41=", $true, $false, $false, $false, $false, $true, 1, $false, "52+19=", 2)
$d.Content.Find.Execute("81-34=", $true, $false, $false, $false, $false, $true, 1, $false, "30-28=", 2)
$d.Content.Find.Execute("16+80=", $true, $false, $false, $false, $false, $true, 1, $false, "95-39=", 2)
$d.Content.Find.Execute("13+32=", $true, $false, $false, $false, $false, $true, 1, $false, "7+11=", 2)
$d.Content.Find.Execute("98-34=", $true, $false, $false, $false, $false, $true, 1, $false, "3+18=", 2)
$d.Content.Find.Execute("74-60=", $true, $false, $false, $false, $false, $true, 1, $false, "22+49=", 2)
$d.Content.Find.Execute("41-6=", $true, $false, $false, $false, $false, $true, 1, $false, "63-28=", 2)
$d.Content.Find.Execute("86-62=", $true, $false, $false, $false, $false, $true, 1, $false, "90-13=", 2)
$d.Content.Find.Execute("49+32=", $true, $false, $false, $false, $false, $true, 1, $false, "80+0=", 2)
$d.Content.Find.Execute("88-61=", $true, $false, $false, $false, $false, $true, 1, $false, "80-16=", 2)
$d.Content.Find.Execute("34+14=", $true, $false, $false, $false, $false, $true, 1, $false, "16+33=", 2)
$d.Content.Find.Execute("22-11=", $true, $false, $false, $false, $false, $true, 1, $false, "11+12=", 2)
$d.Content.Find.Execute("12+20=", $true, $false, $false, $false, $false, $true, 1, $false, "76-40=", 2)
$d.Content.Find.Execute("57+24=", $true, $false, $false, $false, $false, $true, 1, $false, "73+20=", 2)
$d.Content.Find.Execute("42+15=", $true, $false, $false, $false, $false, $true, 1, $false, "82-8=", 2)
$d.Content.Find.Execute("59-58=", $true, $false, $false, $false, $false, $true, 1, $false, "22-5=", 2)
$d.Content.Find.Execute("90-52=", $true, $false, $false, $false, $false, $true, 1, $false, "83-13=", 2)
$d.Content.Find.Execute("52-13=", $true, $false, $false, $false, $false, $true, 1, $false, "23+9=", 2)
$d.Content.Find.Execute("35+33=", $true, $false, $false, $false, $false, $true, 1, $false, "40+57=", 2)
$d.Content.Find.Execute("6+38=", $true, $false, $false, $false, $false, $true, 1, $false, "11+81=", 2)
$d.Content.Find.Execute("30+68=", $true, $false, $false, $false, $false, $true, 1, $false, "30+50=", 2)
$d.Content.Find.Execute("14+17=", $true, $false, $false, $false, $false, $true, 1, $false, "62+18=", 2)
$d.Content.Find.Execute("98-2=", $true, $false, $false, $false, $false, $true, 1, $false, "40+50=", 2)
$d.Content.Find.Execute("12+65=", $true, $false, $false, $false, $false, $true, 1, $false, "36+55=", 2)
$d.Content.Find.Execute("21+70=", $true, $false, $false, $false, $false, $true, 1, $false, "28-9=", 2)
$d.Content.Find.Execute("45+40=", $true, $false, $false, $false, $false, $true, 1, $false, "24+46=", 2)
$d.Content.Find.Execute("61-43=", $true, $false, $false, $false, $false, $true, 1, $false, "10-0=", 2)
$d.Content.Find.Execute("26+45=", $true, $false, $false, $false, $false, $true, 1, $false, "43-39=", 2)
$d.Content.Find.Execute("67+12=", $true, $false, $false, $false, $false, $true, 1, $false, "41-22=", 2)
$d.Content.Find.Execute("69-16=", $true, $false, $false, $false, $false, $true, 1, $false, "75-0=", 2)
$d.Content.Find.Execute("18+34=", $true, $false, $false, $false, $false, $true, 1, $false, "37+55=", 2)
$d.Content.Find.Execute("60+39=", $true, $false, $false, $false, $false, $true, 1, $false, "54-52=", 2)
$d.Content.Find.Execute("8+73=", $true, $false, $false, $false, $false, $true, 1, $false, "98-97=", 2)
$d.Content.Find.Execute("65+1=", $true, $false, $false, $false, $false, $true, 1, $false, "74+10=", 2)
$d.Content.Find.Execute("31+49=", $true, $false, $false, $false, $false, $true, 1, $false, "34-9=", 2)
